$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.425.47"
$ws.Range("D3").Value = "1.802.06"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Formula = "'228.25"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").Formula = "'0.581"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Formula = "'34.79"
$ws.Range("E8").Value = "  +5.34%  "
$ws.Range("E9").Value = "  +1.32%  "
$ws.Range("D10").Formula = "'0.0695"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("E11").Value = "  +0.21%  "
$ws.Range("D12").Value = "2.063.57"
$ws.Range("E12").Value = "  +0.39%  "
$ws.Range("D13").Formula = "'11.23"
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "1.816.91"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("E15").Value = "  +1.02%  "
$ws.Range("D16").Value = "34.427.84"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Formula = "'69.07"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "0.0₃0799"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("D20").Formula = "'245.54"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("D21").Formula = "'11.51"
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").Formula = "'4.17"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Formula = "'173.44"
$ws.Range("E24").Value = "  +3.85%  "
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Formula = "'7.81"
$ws.Range("E26").Value = "  +6.71%  "
$ws.Range("D27").Formula = "'16.78"
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  +2.17%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Formula = "'0.0531"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Formula = "'3.84"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("E34").Value = "  -0.31%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Formula = "'0.682"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "1.394.34"
$ws.Range("E36").Value = "  -2.41%  "
$ws.Range("E37").Value = "  -4.12%  "
$ws.Range("E38").Value = "  -0.23%  "
$ws.Range("E39").Value = "  -1.31%  "
$ws.Range("D40").Formula = "'83.49"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("E41").Value = "  +2.46%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Formula = "'2.42"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Formula = "'13.58"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").Formula = "'1.11"
$ws.Range("E45").Value = "  +3.49%  "
$ws.Range("D46").Formula = "'0.0511"
$ws.Range("E46").Value = "  -3.51%  "
$ws.Range("E47").Value = "  -2.38%  "
$ws.Range("D48").Value = "1.963.11"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").Formula = "'105.08"
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").Value = "0.0₆0130"
$ws.Range("E51").Value = "  +1.41%  "
